# The post "「草を食む幸運な子羊たち」" (row 723) was removed from posts.xlsx.
# Deleting the entire row shifts every following row up by one, which is
# exactly the change described by the diff (dimension A1:C778 -> A1:C777,
# and rows 724-778 renumbered to 723-777 with unchanged content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("723:723").Delete()
